$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revised daily demand (D) and day-of-week share ratio (F) values per 20201023 update
$ws.Cells.Item(2, 4).Value = 91.76000000000001
$ws.Cells.Item(2, 6).Value = 0.93
$ws.Cells.Item(3, 4).Value = 91.76000000000001
$ws.Cells.Item(3, 6).Value = 1.133
$ws.Cells.Item(4, 4).Value = 91.76000000000001
$ws.Cells.Item(4, 6).Value = 0.903
$ws.Cells.Item(5, 4).Value = 91.76000000000001
$ws.Cells.Item(5, 6).Value = 1.05
$ws.Cells.Item(6, 4).Value = 91.76000000000001
$ws.Cells.Item(6, 6).Value = 0.984
$ws.Cells.Item(7, 4).Value = 98.26300000000001
$ws.Cells.Item(7, 6).Value = 0.864
$ws.Cells.Item(8, 4).Value = 98.26300000000001
$ws.Cells.Item(8, 6).Value = 0.8100000000000001
$ws.Cells.Item(9, 4).Value = 98.26300000000001
$ws.Cells.Item(9, 6).Value = 0.945
$ws.Cells.Item(10, 4).Value = 98.26300000000001
$ws.Cells.Item(10, 6).Value = 1.127
$ws.Cells.Item(11, 4).Value = 98.26300000000001
$ws.Cells.Item(11, 6).Value = 1.254
$ws.Cells.Item(12, 4).Value = 118.55
$ws.Cells.Item(12, 6).Value = 1.241
$ws.Cells.Item(13, 4).Value = 118.55
$ws.Cells.Item(13, 6).Value = 0.963
$ws.Cells.Item(14, 4).Value = 118.55
$ws.Cells.Item(14, 6).Value = 0.769
$ws.Cells.Item(15, 4).Value = 118.55
$ws.Cells.Item(15, 6).Value = 1.027
$ws.Cells.Item(16, 4).Value = 71.833
$ws.Cells.Item(16, 6).Value = 0.974
$ws.Cells.Item(17, 4).Value = 71.833
$ws.Cells.Item(17, 6).Value = 1.104
$ws.Cells.Item(18, 4).Value = 71.833
$ws.Cells.Item(18, 6).Value = 1.081
$ws.Cells.Item(19, 4).Value = 71.833
$ws.Cells.Item(19, 6).Value = 0.841
$ws.Cells.Item(20, 4).Value = 49.04
$ws.Cells.Item(20, 6).Value = 1.073
$ws.Cells.Item(21, 4).Value = 49.04
$ws.Cells.Item(21, 6).Value = 1.134
$ws.Cells.Item(22, 4).Value = 49.04
$ws.Cells.Item(22, 6).Value = 1.149
$ws.Cells.Item(23, 4).Value = 49.04
$ws.Cells.Item(23, 6).Value = 0.849
$ws.Cells.Item(24, 4).Value = 49.04
$ws.Cells.Item(24, 6).Value = 0.795
$ws.Cells.Item(25, 4).Value = 33.387
$ws.Cells.Item(25, 6).Value = 1.292
$ws.Cells.Item(26, 4).Value = 33.387
$ws.Cells.Item(26, 6).Value = 1.075
$ws.Cells.Item(27, 4).Value = 33.387
$ws.Cells.Item(27, 6).Value = 0.838
$ws.Cells.Item(28, 4).Value = 33.387
$ws.Cells.Item(28, 6).Value = 0.795
$ws.Cells.Item(29, 4).Value = 20.408
$ws.Cells.Item(29, 6).Value = 1.172
$ws.Cells.Item(30, 4).Value = 20.408
$ws.Cells.Item(30, 6).Value = 1.013
$ws.Cells.Item(31, 4).Value = 20.408
$ws.Cells.Item(31, 6).Value = 0.9399999999999999
$ws.Cells.Item(32, 4).Value = 20.408
$ws.Cells.Item(32, 6).Value = 0.875
$ws.Cells.Item(33, 4).Value = 18.752
$ws.Cells.Item(33, 6).Value = 0.963
$ws.Cells.Item(34, 4).Value = 18.752
$ws.Cells.Item(34, 6).Value = 1.068
$ws.Cells.Item(35, 4).Value = 18.752
$ws.Cells.Item(35, 6).Value = 1.148
$ws.Cells.Item(36, 4).Value = 18.752
$ws.Cells.Item(36, 6).Value = 0.912
$ws.Cells.Item(37, 4).Value = 18.752
$ws.Cells.Item(37, 6).Value = 0.908
$ws.Cells.Item(38, 4).Value = 24.203
$ws.Cells.Item(38, 6).Value = 0.666
$ws.Cells.Item(39, 4).Value = 24.203
$ws.Cells.Item(39, 6).Value = 0.6870000000000001
$ws.Cells.Item(40, 4).Value = 24.203
$ws.Cells.Item(40, 6).Value = 0.905
$ws.Cells.Item(41, 4).Value = 24.203
$ws.Cells.Item(41, 6).Value = 1.743
$ws.Cells.Item(42, 4).Value = 62.442
$ws.Cells.Item(42, 6).Value = 0.694
$ws.Cells.Item(43, 4).Value = 62.442
$ws.Cells.Item(43, 6).Value = 0.856
$ws.Cells.Item(44, 4).Value = 62.442
$ws.Cells.Item(44, 6).Value = 0.9360000000000001
$ws.Cells.Item(45, 4).Value = 62.442
$ws.Cells.Item(45, 6).Value = 1.05
$ws.Cells.Item(46, 4).Value = 62.442
$ws.Cells.Item(46, 6).Value = 1.463
$ws.Cells.Item(47, 4).Value = 83.916
$ws.Cells.Item(47, 6).Value = 1.005
$ws.Cells.Item(48, 4).Value = 83.916
$ws.Cells.Item(48, 6).Value = 1.131
$ws.Cells.Item(49, 4).Value = 83.916
$ws.Cells.Item(49, 6).Value = 0.871
$ws.Cells.Item(50, 4).Value = 83.916
$ws.Cells.Item(50, 6).Value = 0.993
$ws.Cells.Item(51, 4).Value = 123.594
$ws.Cells.Item(51, 6).Value = 1.291
$ws.Cells.Item(52, 4).Value = 123.594
$ws.Cells.Item(52, 6).Value = 1.056
$ws.Cells.Item(53, 4).Value = 123.594
$ws.Cells.Item(53, 6).Value = 0.793
$ws.Cells.Item(54, 4).Value = 123.594
$ws.Cells.Item(54, 6).Value = 0.859
